$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values per row for columns D, J, K, L, M, O, P (weekly shift of the dataset).
# Row 3 is left untouched; rows 2 and 4-16 get the rotated values below.
$data = @{
    2  = @{ D = 44224; J = 80;  K = 8500; L = 9000;  M = 8719; O = "Región del Maule";             P = 145 }
    4  = @{ D = 44159; J = 35;  K = 7500; L = 8000;  M = 7714; O = "Región de Arica y Parinacota";  P = 129 }
    5  = @{ D = 44202; J = 50;  K = 8000; L = 9000;  M = 8400; O = "Región del Maule";             P = 140 }
    6  = @{ D = 44271; J = 55;  K = 9000; L = 9500;  M = 9227; O = "Región del Maule";             P = 154 }
    7  = @{ D = 44259; J = 70;  K = 9000; L = 9500;  M = 9214; O = "Región del Maule";             P = 154 }
    8  = @{ D = 44204; J = 45;  K = 9500; L = 10000; M = 9722; O = "Región del Maule";             P = 162 }
    9  = @{ D = 44160; J = 90;  K = 7500; L = 8000;  M = 7667; O = "Región de Arica y Parinacota";  P = 128 }
    10 = @{ D = 44218; J = 65;  K = 9000; L = 10000; M = 9615; O = "Región del Maule";             P = 160 }
    11 = @{ D = 44210; J = 60;  K = 8000; L = 9000;  M = 8417; O = "Región de Arica y Parinacota";  P = 140 }
    12 = @{ D = 44208; J = 100; K = 7000; L = 8000;  M = 7350; O = "Región del Maule";             P = 122 }
    13 = @{ D = 44216; J = 55;  K = 9500; L = 10000; M = 9773; O = "Región del Maule";             P = 163 }
    14 = @{ D = 44162; J = 43;  K = 8000; L = 8500;  M = 8209; O = "Región de Arica y Parinacota";  P = 137 }
    15 = @{ D = 44266; J = 60;  K = 9000; L = 9500;  M = 9208; O = "Región del Maule";             P = 153 }
    16 = @{ D = 44264; J = 43;  K = 8500; L = 9000;  M = 8709; O = "Región del Maule";             P = 145 }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]

    $ws.Cells.Item($row, 4).Value  = $vals.D   # D - Fecha
    $ws.Cells.Item($row, 10).Value = $vals.J   # J - Volumen
    $ws.Cells.Item($row, 11).Value = $vals.K   # K - Precio mínimo
    $ws.Cells.Item($row, 12).Value = $vals.L   # L - Precio máximo
    $ws.Cells.Item($row, 13).Value = $vals.M   # M - Precio promedio ponderado
    $ws.Cells.Item($row, 15).Value = $vals.O   # O - Origen
    $ws.Cells.Item($row, 16).Value = $vals.P   # P - Precio $/Kg
}
